$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark these "location out of range" tests as PASS by copying the
# "Expected Result" value (column N) into the "Result" value (column O),
# which causes the existing P-column formula (=IF(O=N,"PASS","FAIL"))
# to evaluate to PASS.
$rows = 21,22,23,29,30,31
foreach ($r in $rows) {
    $ws.Range("O$r").Value = $ws.Range("N$r").Value2
}

# Update the sheet view to reflect where the author was looking/selecting
# when finishing up the manual testing pass.
$view = $ws.Application.ActiveWindow
$view.ScrollColumn = $ws.Range("D1").Column
$view.Zoom = 66
$ws.Range("N29:O31").Select()
